$d = $word.ActiveDocument
$BR = [char]11

# ---------------------------------------------------------------------------
# This edit performs a content rotation across several paragraphs: the text
# that used to live in one paragraph is relocated to another paragraph,
# while paragraph styles/structure stay put. We therefore:
#   1) capture every "moving" piece of text (with manual line breaks kept as
#      Chr(11), matching how Word represents <w:br/> in Range.Text) from its
#      ORIGINAL location, before making any edits;
#   2) only then write the captured text into its new home.
# This two-phase read-then-write approach avoids clobbering a value before
# it has been copied elsewhere, since several of the moves form a rotation
# cycle.
# ---------------------------------------------------------------------------

function Strip-Para([string]$s) {
    if ($s.Length -gt 0 -and [int][char]$s[$s.Length - 1] -eq 13) {
        return $s.Substring(0, $s.Length - 1)
    }
    return $s
}

# ---- Phase 1: capture original text of every paragraph that changes ------

$txt_P6  = Strip-Para $d.Paragraphs(6).Range.Text    # "Oferecer uma base sólida..."
$txt_P7  = Strip-Para $d.Paragraphs(7).Range.Text    # "Provide a solid foundation..."
$txt_P9  = Strip-Para $d.Paragraphs(9).Range.Text    # "6270264 - Juan Fernando Zapata Zapata"
$txt_P11 = Strip-Para $d.Paragraphs(11).Range.Text   # "Números reais, funções reais..."
$txt_P12 = Strip-Para $d.Paragraphs(12).Range.Text   # "Real numbers, real functions..."
$txt_P14 = Strip-Para $d.Paragraphs(14).Range.Text   # "Funções Reais: ...<br>Matrizes...<br>Modelagem..."
$txt_P19 = Strip-Para $d.Paragraphs(19).Range.Text   # "Leithold...<br><br>ANTON...<br><br>THOMAS...<br><br>FLEMMING..."

# Paragraph 17 holds three bold labels ("Método: ", "Critério: ",
# "Norma de recuperação: "), each followed by a plain-text value. Capture the
# three values (with their trailing Chr(11) break, when present) using the
# labels as stable anchors.
$p17 = $d.Paragraphs(17)

$lblMetodo = $p17.Range.Duplicate
[void]$lblMetodo.Find.Execute("M" + [char]233 + "todo: ")

$lblCriterio = $p17.Range.Duplicate
[void]$lblCriterio.Find.Execute("Crit" + [char]233 + "rio: ")

$lblNorma = $p17.Range.Duplicate
[void]$lblNorma.Find.Execute("Norma de recupera" + [char]231 + [char]227 + "o: ")

$valMetodoRange   = $d.Range($lblMetodo.End, $lblCriterio.Start)
$valCriterioRange = $d.Range($lblCriterio.End, $lblNorma.Start)
$valNormaRange    = $d.Range($lblNorma.End, $p17.Range.End)

$txt_valMetodo   = $valMetodoRange.Text
$txt_valCriterio = $valCriterioRange.Text
$txt_valNorma    = Strip-Para $valNormaRange.Text

# ---- Phase 2: write captured text into its new home -----------------------

# Cycle A (length 2): P7 <-> P12
$d.Paragraphs(7).Range.Text  = $txt_P12
$d.Paragraphs(12).Range.Text = $txt_P7

# Cycle B (length 8):
#   P6 -> P9 -> P19 -> (P17 Norma-value) -> (P17 Criterio-value)
#      -> (P17 Metodo-value) -> P14 -> P11 -> P6
$d.Paragraphs(9).Range.Text  = $txt_P6
$d.Paragraphs(19).Range.Text = $txt_P9
$d.Paragraphs(14).Range.Text = $txt_valMetodo
$d.Paragraphs(11).Range.Text = $txt_P14
$d.Paragraphs(6).Range.Text  = $txt_P11

# Re-locate paragraph 17's labels (their own text/position is unchanged) and
# rewrite the three values, starting with the right-most one so that the
# offsets used to find the labels further left stay valid while we still
# need them.
$lblNorma2 = $p17.Range.Duplicate
[void]$lblNorma2.Find.Execute("Norma de recupera" + [char]231 + [char]227 + "o: ")
$valNormaRange2 = $d.Range($lblNorma2.End, $p17.Range.End)
$valNormaRange2.Text = $txt_P19

$lblCriterio2 = $p17.Range.Duplicate
[void]$lblCriterio2.Find.Execute("Crit" + [char]233 + "rio: ")
$lblNorma3 = $p17.Range.Duplicate
[void]$lblNorma3.Find.Execute("Norma de recupera" + [char]231 + [char]227 + "o: ")
$valCriterioRange2 = $d.Range($lblCriterio2.End, $lblNorma3.Start)
$valCriterioRange2.Text = $txt_valNorma

$lblMetodo2 = $p17.Range.Duplicate
[void]$lblMetodo2.Find.Execute("M" + [char]233 + "todo: ")
$lblCriterio3 = $p17.Range.Duplicate
[void]$lblCriterio3.Find.Execute("Crit" + [char]233 + "rio: ")
$valMetodoRange2 = $d.Range($lblMetodo2.End, $lblCriterio3.Start)
$valMetodoRange2.Text = $txt_valCriterio
